$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2039.4117  # H17
$ws.Cells.Item(17, 10).Value = 2191.3333  # J17
$ws.Cells.Item(17, 12).Value = 6573.999899999999  # L17
$ws.Cells.Item(17, 14).Value = -6909.999899999999  # N17

$ws.Cells.Item(40, 8).Value = 1757.7037  # H40
$ws.Cells.Item(40, 9).Value = 2277.9  # I40
$ws.Cells.Item(40, 10).Value = 1451.7059  # J40
$ws.Cells.Item(40, 11).Value = 2277.9  # K40
$ws.Cells.Item(40, 12).Value = 1451.7059  # L40
$ws.Cells.Item(40, 13).Value = -2102.9  # M40
$ws.Cells.Item(40, 14).Value = -1801.7059  # N40

$ws.Cells.Item(58, 8).Value = 2398.8813  # H58
$ws.Cells.Item(58, 9).Value = 353.4  # I58
$ws.Cells.Item(58, 10).Value = 2816.3264  # J58
$ws.Cells.Item(58, 11).Value = 1060.2  # K58
$ws.Cells.Item(58, 12).Value = 8448.9792  # L58
$ws.Cells.Item(58, 13).Value = -910.1999999999998  # M58
$ws.Cells.Item(58, 14).Value = -8748.9792  # N58

$ws.Cells.Item(69, 8).Value = 10420001  # H69
$ws.Cells.Item(69, 9).Value = 11366955  # I69
$ws.Cells.Item(69, 11).Value = 34100865  # K69
$ws.Cells.Item(69, 13).Value = -34099991  # M69

$ws.Cells.Item(72, 8).Value = 10420001  # H72
$ws.Cells.Item(72, 9).Value = 11366955  # I72
$ws.Cells.Item(72, 11).Value = 102302595  # K72
$ws.Cells.Item(72, 13).Value = -102298227  # M72

$ws.Cells.Item(76, 8).Value = 52815.15  # H76
$ws.Cells.Item(76, 9).Value = 55410.684  # I76
$ws.Cells.Item(76, 11).Value = 55410.684  # K76
$ws.Cells.Item(76, 13).Value = -55095.684  # M76

$ws.Cells.Item(79, 8).Value = 52815.15  # H79
$ws.Cells.Item(79, 9).Value = 55410.684  # I79
$ws.Cells.Item(79, 11).Value = 55410.684  # K79
$ws.Cells.Item(79, 13).Value = -54318.684  # M79

$ws.Cells.Item(86, 8).Value = 83418320  # H86
$ws.Cells.Item(86, 9).Value = 112888.336  # I86
$ws.Cells.Item(86, 10).Value = 333334600  # J86
$ws.Cells.Item(86, 11).Value = 112888.336  # K86
$ws.Cells.Item(86, 12).Value = 333334600  # L86
$ws.Cells.Item(86, 13).Value = -111765.336  # M86
$ws.Cells.Item(86, 14).Value = -333336846  # N86

$ws.Cells.Item(89, 8).Value = 83418320  # H89
$ws.Cells.Item(89, 9).Value = 112888.336  # I89
$ws.Cells.Item(89, 10).Value = 333334600  # J89
$ws.Cells.Item(89, 11).Value = 564441.6799999999  # K89
$ws.Cells.Item(89, 12).Value = 1666673000  # L89
$ws.Cells.Item(89, 13).Value = -558825.6799999999  # M89
$ws.Cells.Item(89, 14).Value = -1666684232  # N89

$ws.Cells.Item(103, 8).Value = 348.66666  # H103
$ws.Cells.Item(103, 9).Value = 348.66666  # I103
$ws.Cells.Item(103, 11).Value = 1045.99998  # K103
$ws.Cells.Item(103, 13).Value = -459.9999800000001  # M103

$ws.Cells.Item(113, 8).Value = 3465.5557  # H113
$ws.Cells.Item(113, 9).Value = 3333.3333  # I113
$ws.Cells.Item(113, 10).Value = 3531.6667  # J113
$ws.Cells.Item(113, 11).Value = 3333.3333  # K113
$ws.Cells.Item(113, 12).Value = 3531.6667  # L113
$ws.Cells.Item(113, 13).Value = -79.33329999999978  # M113
$ws.Cells.Item(113, 14).Value = -10039.6667  # N113

$ws.Cells.Item(125, 8).Value = 125001980  # H125
$ws.Cells.Item(125, 10).Value = 2322.6667  # J125
$ws.Cells.Item(125, 12).Value = 20904.0003  # L125
$ws.Cells.Item(125, 14).Value = -25824.0003  # N125

$ws.Cells.Item(129, 8).Value = 929636.2  # H129
$ws.Cells.Item(129, 9).Value = 546  # I129
$ws.Cells.Item(129, 10).Value = 1115454.2  # J129
$ws.Cells.Item(129, 11).Value = 1638  # K129
$ws.Cells.Item(129, 12).Value = 3346362.6  # L129
$ws.Cells.Item(129, 13).Value = 3362  # M129
$ws.Cells.Item(129, 14).Value = -3356362.6  # N129

$ws.Cells.Item(141, 8).Value = 3285.818  # H141
$ws.Cells.Item(141, 9).Value = 2588.2222  # I141
$ws.Cells.Item(141, 10).Value = 6425  # J141
$ws.Cells.Item(141, 11).Value = 7764.6666  # K141
$ws.Cells.Item(141, 12).Value = 19275  # L141
$ws.Cells.Item(141, 13).Value = -2584.6666  # M141
$ws.Cells.Item(141, 14).Value = -29635  # N141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(135, 8).Value = 45230.855  # H135
$ws.Cells.Item(135, 10).Value = 45230.855  # J135
$ws.Cells.Item(135, 12).Value = 45230.855  # L135
$ws.Cells.Item(135, 14).Value = -55370.855  # N135

$ws.Cells.Item(139, 8).Value = 48750  # H139
$ws.Cells.Item(139, 10).Value = 51666.668  # J139
$ws.Cells.Item(139, 12).Value = 51666.668  # L139
$ws.Cells.Item(139, 14).Value = -61946.668  # N139

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(81, 8).Value = 8894.286  # H81
$ws.Cells.Item(81, 10).Value = 8894.286  # J81
$ws.Cells.Item(81, 12).Value = 8894.286  # L81
$ws.Cells.Item(81, 14).Value = -11016.286  # N81

$ws.Cells.Item(84, 8).Value = 8894.286  # H84
$ws.Cells.Item(84, 10).Value = 8894.286  # J84
$ws.Cells.Item(84, 12).Value = 26682.858  # L84
$ws.Cells.Item(84, 14).Value = -37290.858  # N84

$ws.Cells.Item(132, 8).Value = 0  # H132
$ws.Cells.Item(132, 10).Value = 0  # J132
$ws.Cells.Item(132, 12).Value = 0  # L132
$ws.Cells.Item(132, 14).ClearContents()  # N132

$ws.Cells.Item(141, 8).Value = 80000  # H141
$ws.Cells.Item(141, 10).Value = 80000  # J141
$ws.Cells.Item(141, 12).Value = 80000  # L141
$ws.Cells.Item(141, 14).Value = -90360  # N141

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 41670016  # H62
$ws.Cells.Item(62, 9).Value = 3359.6  # I62
$ws.Cells.Item(62, 10).Value = 111114450  # J62
$ws.Cells.Item(62, 11).Value = 3359.6  # K62
$ws.Cells.Item(62, 12).Value = 111114450  # L62
$ws.Cells.Item(62, 13).Value = -2735.6  # M62
$ws.Cells.Item(62, 14).Value = -111115698  # N62

$ws.Cells.Item(65, 8).Value = 41670016  # H65
$ws.Cells.Item(65, 9).Value = 3359.6  # I65
$ws.Cells.Item(65, 10).Value = 111114450  # J65
$ws.Cells.Item(65, 11).Value = 16798  # K65
$ws.Cells.Item(65, 12).Value = 555572250  # L65
$ws.Cells.Item(65, 13).Value = -13678  # M65
$ws.Cells.Item(65, 14).Value = -555578490  # N65

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 834.5700000000001  # H131
$ws.Cells.Item(131, 9).Value = 1150  # I131
$ws.Cells.Item(131, 10).Value = 817.96844  # J131
$ws.Cells.Item(131, 11).Value = 3450  # K131
$ws.Cells.Item(131, 12).Value = 2453.90532  # L131
$ws.Cells.Item(131, 13).Value = 1590  # M131
$ws.Cells.Item(131, 14).Value = -12533.90532  # N131

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 5874.483  # H132
$ws.Cells.Item(132, 9).Value = 6319.522  # I132
$ws.Cells.Item(132, 11).Value = 18958.566  # K132
$ws.Cells.Item(132, 13).Value = -16428.566  # M132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 687.25  # H22
$ws.Cells.Item(22, 9).Value = 1034.6666  # I22
$ws.Cells.Item(22, 10).Value = 478.8  # J22
$ws.Cells.Item(22, 11).Value = 1034.6666  # K22
$ws.Cells.Item(22, 12).Value = 478.8  # L22
$ws.Cells.Item(22, 13).Value = -739.6666  # M22
$ws.Cells.Item(22, 14).Value = -1068.8  # N22

$ws.Cells.Item(27, 8).Value = 687.25  # H27
$ws.Cells.Item(27, 9).Value = 1034.6666  # I27
$ws.Cells.Item(27, 10).Value = 478.8  # J27
$ws.Cells.Item(27, 11).Value = 1034.6666  # K27
$ws.Cells.Item(27, 12).Value = 478.8  # L27
$ws.Cells.Item(27, 13).Value = -927.6666  # M27
$ws.Cells.Item(27, 14).Value = -692.8  # N27

$ws.Cells.Item(68, 8).Value = 15626168  # H68
$ws.Cells.Item(68, 9).Value = 1223.3334  # I68
$ws.Cells.Item(68, 11).Value = 1223.3334  # K68
$ws.Cells.Item(68, 13).Value = -474.3334  # M68

$ws.Cells.Item(71, 8).Value = 15626168  # H71
$ws.Cells.Item(71, 9).Value = 1223.3334  # I71
$ws.Cells.Item(71, 11).Value = 6116.666999999999  # K71
$ws.Cells.Item(71, 13).Value = -2372.666999999999  # M71

$ws.Cells.Item(132, 8).Value = 8747.143  # H132
$ws.Cells.Item(132, 9).Value = 10838.6  # I132
$ws.Cells.Item(132, 10).Value = 3518.5  # J132
$ws.Cells.Item(132, 11).Value = 32515.8  # K132
$ws.Cells.Item(132, 12).Value = 10555.5  # L132
$ws.Cells.Item(132, 13).Value = -29985.8  # M132
$ws.Cells.Item(132, 14).Value = -15615.5  # N132

$ws.Cells.Item(136, 8).Value = 8072  # H136
$ws.Cells.Item(136, 9).Value = 9762.666999999999  # I136
$ws.Cells.Item(136, 11).Value = 29288.001  # K136
$ws.Cells.Item(136, 13).Value = -26738.001  # M136
